$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Reposition the selection (cosmetic, matches author's later cursor spot) ---
$excel.ActiveWindow.ScrollRow = 31
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("G42").Select()

# --- Replace the hyperlinked "(Author et al., year)" references with plain
#     text "Author et al., year" (no parens, no hyperlink). The source cells
#     currently hold shared-string text "(Rammstedt et al., 2020)" for
#     E17:E31 and "(Kovaleva et al., 2014)" for E32:E35, wrapped in
#     hyperlinks that point at a (non-existent) Sources sheet. ---
$ws.Range("E17:E31").Value = "Rammstedt et al., 2020"
$ws.Range("E32:E35").Value = "Kovaleva et al., 2014"

# Drop the hyperlink objects entirely.
$ws.Hyperlinks.Delete()

# Reset the E-column cells back to the plain default style (no border / no
# fill / no special "Link" font) -- copy the already-plain style used by B1.
$plainStyle = $ws.Range("B1").Style
$ws.Range("E17:E35").Style = $plainStyle

# The "Link" cell style (used only by the hyperlinked cells above) is now
# unused -- remove it so its backing font/fill/border entries are dropped too.
$wb.Styles("Link").Delete()
